$d = $word.ActiveDocument

# 1. Merge "Resume Date: " and "02/20/2023" runs into one
$d.Content.Find.Execute("Resume Date: 02/20/2023", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Resume Date: 02/20/2023", 2)

# 2. Merge the github hyperlink runs into one
$d.Content.Find.Execute("https://github.com/lee-lindley", $true, $false, $false, $false, $false,
                         $true, 1, $false, "https://github.com/lee-lindley", 2)

# 3. Split "work in process." into three runs: "work in pro", "gr", "ess."
$d.Content.Find.Execute("work in process.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "work in progress.", 2)
